$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update marking scheme (row 11): right mark 4 -> 5, wrong mark -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update totals (row 12) derived from new marking scheme:
# Right total = 22 * 5 = 110 ; Wrong total stays at -0 (0 wrong answers)
$ws.Range("B12").Value = 110
$ws.Range("C12").Value = -0

# Update the fraction display: new total / new max (28 * 5 = 140)
$ws.Range("E12").Value = "110.0/140"
